# The report's "Age" checklist is repeated three times in the document
# ("☐ Young ☐ Adult ☐ Old"). In the middle occurrence, the word "Young "
# and the single space that follows it live in two separate (but
# identically formatted) runs:
#
#     <w:r><w:t xml:space="preserve">Young </w:t></w:r>
#     <w:r><w:t xml:space="preserve"> </w:t></w:r>
#
# The edit collapses that redundant run split into one run:
#
#     <w:r><w:t xml:space="preserve">Young  </w:t></w:r>
#
# Locate the unique, unambiguous anchor text "Young   Adult" (three
# spaces - the trailing space of "Young ", the lone-space run, and the
# space hidden behind the Wingdings checkbox symbol that precedes
# "Adult") so we edit only the intended paragraph and leave the other
# two "Young"/"Adult"/"Old" checklists untouched.

$d = $word.ActiveDocument

$anchor = $d.Content
$found = $anchor.Find.Execute("Young   Adult", $true, $false, $false, $false, `
                               $false, $true, 0, $false, "", 0)

if ($found) {
    # Expand from the matched text out to the whole paragraph, then grab a
    # fresh Range over those bounds so a later Find starts cleanly at the
    # paragraph start instead of resuming from the previous match.
    $paraRange = $d.Range($anchor.Start, $anchor.End)
    $null = $paraRange.Expand(4)

    $editRange = $d.Range($paraRange.Start, $paraRange.End)

    # Replace "Young" + the two adjoining single-space runs with a single
    # run containing "Young" followed by two spaces - same visible text,
    # merged into one <w:r>, scoped to this paragraph only (Wrap = 0 /
    # wdFindStop keeps the search from drifting into the next checklist).
    $null = $editRange.Find.Execute("Young  ", $true, $false, $false, $false, `
                                     $false, $true, 0, $false, "Young  ", 2)
}
